# Update the simulation-result grid on Sheet1 (rows 2-25, columns B:O)
# with the refreshed values for the "380 kV" case.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is: row number, then the 14 new values for columns B..O
$data = @(
    ,(2, 0, 0.01681254295004209, 0.01899188952842934, 0.1146634113283724, 0.6091378076441956, 0.00240415971103272, 0, 0.4538774918128325, 0, 0.8380131953018974, 0, 0.3298747082999256, 1.07339142281571, 2.024516934344973)
    ,(3, 0, 0.0147254155090053, 0.01751924932826654, 0.1086312128689926, 0.6032880531546496, 0.002406895151447619, 0, 0.4502697176569157, 0, 0.7336431570084301, 0, 0.2950121864063036, 1.080729859088088, 2.018984902989047)
    ,(4, 0, 0.01343831474956403, 0.0166085100554838, 0.1050340148368463, 0.6001579022374557, 0.002408664694095223, 0, 0.4484003038367916, 0, 0.669418434303708, 0, 0.2736761886250108, 1.085675331145673, 2.017113887164697)
    ,(5, 0, 0.01291242623318567, 0.01623575885282236, 0.1035946732239452, 0.5989982290491795, 0.002409408489591541, 0, 0.4477252976266541, 0, 0.6432122760965058, 0, 0.2649991893264598, 1.087801311374797, 2.016734348431726)
    ,(6, 0, 0.01282501989250306, 0.0161737668409927, 0.1033572673889651, 0.5988126600442953, 0.002409533368718176, 0, 0.4476184511092924, 0, 0.6388587480679462, 0, 0.2635594445772611, 1.088161015673265, 2.01669443218799)
    ,(7, 0, 0.01343122800886221, 0.01660348951912027, 0.1050144962311848, 0.6001417935033828, 0.002408674633392417, 0, 0.4483908492402904, 0, 0.6690651449230813, 0, 0.2735590962646697, 1.085703554654181, 2.017107219132726)
    ,(8, 0, 0.01609407586656886, 0.01848549352083495, 0.1125612274958598, 0.6070248847405111, 0.002405084258752954, 0, 0.4525616715288052, 0, 0.8020565086169427, 0, 0.3178396617592227, 1.075830531046172, 2.022292390320246)
    ,(9, 0, 0.0212707540093362, 0.02212333300170854, 0.1282186494446265, 0.624195624161338, 0.002398754302800782, 0, 0.4634922526956373, 0, 1.061682877641772, 0, 0.4052290159973779, 1.059953605942844, 2.044602706877185)
    ,(10, 0, 0.02504586729106961, 0.02476287448612879, 0.1402639348223857, 0.6390667673136221, 0.002394532604013455, 0, 0.4732132270109872, 0, 1.251669171616072, 0, 0.4697819836003987, 1.05040764037696, 2.068452627410295)
    ,(11, 0, 0.02675702266263613, 0.02595626656133021, 0.1458653543728019, 0.6463255797945919, 0.002392704248664137, 0, 0.4780054338546975, 0, 1.337924491027138, 0, 0.4992268191999614, 1.046524071712767, 2.080934581418376)
    ,(12, 0, 0.02740408803884975, 0.02640709421983445, 0.1480043068697228, 0.6491455682997156, 0.002392025073538084, 0, 0.4798735320323857, 0, 1.370561454186941, 0, 0.5103882744894861, 1.045119392641936, 2.085896824399754)
    ,(13, 0, 0.0272647719068857, 0.0263100490622179, 0.1475428483363999, 0.6485350617164869, 0.002392170760687976, 0, 0.4794688261900362, 0, 1.363533683634671, 0, 0.5079839524862138, 1.045418983277465, 2.084817624856299)
    ,(14, 0, 0.0268102755470494, 0.02599337830627491, 0.1460409686126596, 0.6465561530700938, 0.002392648108565583, 0, 0.4781580523292135, 0, 1.340610084662558, 0, 0.5001448521141896, 1.046407186572694, 2.081338101563802)
    ,(15, 0, 0.02653176386233724, 0.02579926634326313, 0.1451233524394766, 0.6453532975767757, 0.002392942212847665, 0, 0.4773621240906962, 0, 1.32656527939821, 0, 0.4953446506573869, 1.047021076084881, 2.079237499949983)
    ,(16, 0, 0.02493391304534498, 0.02468473329953014, 0.1399003468335565, 0.6386023431281842, 0.002394653939720448, 0, 0.472907507125484, 0, 1.246028622059782, 0, 0.4678592847410386, 1.050670671275036, 2.067669824860303)
    ,(17, 0, 0.02395208654687053, 0.023999101768581, 0.136727636893859, 0.6345874947461709, 0.002395727577269957, 0, 0.470269657025959, 0, 1.196577242791989, 0, 0.4510181751292208, 1.053027087027544, 2.06099215680041)
    ,(18, 0, 0.02338678740863998, 0.02360405458894377, 0.1349142406759682, 0.6323247311246121, 0.002396353779288327, 0, 0.468787259293272, 0, 1.168118171039509, 0, 0.4413390892936775, 1.054425635203863, 2.057304947416014)
    ,(19, 0, 0.02319528826116368, 0.02347018071642992, 0.1343022161031513, 0.6315665730210611, 0.002396567292309009, 0, 0.4682913205736696, 0, 1.158479726892324, 0, 0.4380632033579275, 1.05490658090099, 2.056082878331154)
    ,(20, 0, 0.02405666371463155, 0.02407216003265944, 0.1370641889273685, 0.6350100710549782, 0.002395612389535645, 0, 0.4705468552583696, 0, 1.201843088175337, 0, 0.4528101669332898, 1.052771771876948, 2.061687101909484)
    ,(21, 0, 0.02694379710345629, 0.02608642182050858, 0.1464816213328106, 0.6471354715122857, 0.002392507542567225, 0, 0.4785416079077933, 0, 1.347344020799596, 0, 0.5024470789004027, 1.046115137982319, 2.082353721108205)
    ,(22, 0, 0.02882537767685278, 0.02739652723627017, 0.1527404358821016, 0.6554754270235321, 0.002390555161461425, 0, 0.484077939802134, 0, 1.442284710790659, 0, 0.5349537044084371, 1.042148996122748, 2.097234187724865)
    ,(23, 0, 0.02782163877198229, 0.02669788761167524, 0.1493903824007248, 0.6509861635437062, 0.002391590175910581, 0, 0.4810945536052742, 0, 1.391627550053215, 0, 0.5175982899231286, 1.04423064761275, 2.089166235800349)
    ,(24, 0, 0.02400938697346078, 0.02403913309494499, 0.1369120006429867, 0.6348188828215484, 0.002395664437921762, 0, 0.4704214276139496, 0, 1.199462489118162, 0, 0.4519999978046627, 1.052887063372452, 2.061372444214385)
    ,(25, 0, 0.01987522836559918, 0.02114495490579316, 0.123889101592674, 0.619155540592736, 0.002400391086528455, 0, 0.4602393351252019, 0, 0.9915772675136907, 0, 0.3815274550801036, 1.063876299924999, 2.03726144413352)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $arr = New-Object 'object[,]' 1,14
    for ($c = 0; $c -lt 14; $c++) {
        $arr[0, $c] = $entry[$c + 1]
    }
    $addr = "B" + $r + ":O" + $r
    $rng = $ws.Range($addr)
    $rng.Value = $arr
}

Write-Host "Updated $($data.Count) rows"